$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10 (Q0-Q8) with revised forecast-error values
$ws.Cells.Item(2,2).Value = -0.03905235774959507
$ws.Cells.Item(2,3).Value = 0.4387070408530346
$ws.Cells.Item(2,4).Value = 0.45142052592237
$ws.Cells.Item(2,5).Value = 0.6718783564919844
$ws.Cells.Item(2,6).Value = 0.6774166724122599
$ws.Cells.Item(2,7).Value = 51

$ws.Cells.Item(3,2).Value = 0.0876619451090646
$ws.Cells.Item(3,3).Value = 0.4671253201736847
$ws.Cells.Item(3,4).Value = 0.4722669860490882
$ws.Cells.Item(3,5).Value = 0.6872168406326261
$ws.Cells.Item(3,6).Value = 0.6885227972073756
$ws.Cells.Item(3,7).Value = 50

$ws.Cells.Item(4,2).Value = -0.02200267238528682
$ws.Cells.Item(4,3).Value = 0.4595796052308628
$ws.Cells.Item(4,4).Value = 0.413274940495816
$ws.Cells.Item(4,5).Value = 0.642864636214978
$ws.Cells.Item(4,6).Value = 0.6491460788714786
$ws.Cells.Item(4,7).Value = 49

$ws.Cells.Item(5,2).Value = 0.1002432192375009
$ws.Cells.Item(5,3).Value = 0.5011266898800307
$ws.Cells.Item(5,4).Value = 0.4866691161701831
$ws.Cells.Item(5,5).Value = 0.6976167401734158
$ws.Cells.Item(5,6).Value = 0.6976827882510105
$ws.Cells.Item(5,7).Value = 48

$ws.Cells.Item(6,2).Value = 0.02274536467644276
$ws.Cells.Item(6,3).Value = 0.4979849786295887
$ws.Cells.Item(6,4).Value = 0.4567703583966128
$ws.Cells.Item(6,5).Value = 0.6758478811068456
$ws.Cells.Item(6,6).Value = 0.6827675668981734
$ws.Cells.Item(6,7).Value = 47

$ws.Cells.Item(7,2).Value = 0.092078373345108
$ws.Cells.Item(7,3).Value = 0.5065592491403523
$ws.Cells.Item(7,4).Value = 0.4958636947812035
$ws.Cells.Item(7,5).Value = 0.7041758976145118
$ws.Cells.Item(7,6).Value = 0.7058442120435615
$ws.Cells.Item(7,7).Value = 46

$ws.Cells.Item(8,2).Value = 0.03822042347955541
$ws.Cells.Item(8,3).Value = 0.5402057006093556
$ws.Cells.Item(8,4).Value = 0.4987337210178167
$ws.Cells.Item(8,5).Value = 0.7062108191027781
$ws.Cells.Item(8,6).Value = 0.7131441491908896
$ws.Cells.Item(8,7).Value = 45

$ws.Cells.Item(9,2).Value = 0.1352597487784972
$ws.Cells.Item(9,3).Value = 0.5335606825057502
$ws.Cells.Item(9,4).Value = 0.501733990090413
$ws.Cells.Item(9,5).Value = 0.70833183614067
$ws.Cells.Item(9,6).Value = 0.7033360171492812
$ws.Cells.Item(9,7).Value = 44

$ws.Cells.Item(10,2).Value = 0.05955804992731925
$ws.Cells.Item(10,3).Value = 0.54934389505087
$ws.Cells.Item(10,4).Value = 0.4945675297751405
$ws.Cells.Item(10,5).Value = 0.7032549536086756
$ws.Cells.Item(10,6).Value = 0.7090213886885923
$ws.Cells.Item(10,7).Value = 43

# Add new row 11 (Q9) with its data, copying the format from row 10's label cell
$ws.Cells.Item(11,1).Value = "Q9"
$ws.Cells.Item(10,1).Copy()
$ws.Cells.Item(11,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(11,2).Value = 0.1423836843877604
$ws.Cells.Item(11,3).Value = 0.5466645995897268
$ws.Cells.Item(11,4).Value = 0.4954308102877077
$ws.Cells.Item(11,5).Value = 0.7038684609269744
$ws.Cells.Item(11,6).Value = 0.6976724939559396
$ws.Cells.Item(11,7).Value = 42
